# Auto-generated Excel COM-interop script
# Applies updated currentAveragePrice / LevePrice / LeveProfit values
# per sheet, matching the target diff for Durandal_Profits.xlsx

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 166.35
$ws.Range("I15").Value = 166.35
$ws.Range("K15").Value = 499.05
$ws.Range("M15").Value = -330.05
$ws.Range("H33").Value = 6821.533
$ws.Range("I33").Value = 91.44444
$ws.Range("K33").Value = 91.44444
$ws.Range("M33").Value = 137.55556
$ws.Range("H137").Value = 1368.6111
$ws.Range("I137").Value = 1202.7693
$ws.Range("K137").Value = 3608.3079
$ws.Range("M137").Value = -1058.3079
$ws.Range("H139").Value = 70276
$ws.Range("J139").Value = 70276
$ws.Range("L139").Value = 70276
$ws.Range("N139").Value = -80556
$ws.Range("H140").Value = 64342.31
$ws.Range("J140").Value = 88494.44500000001
$ws.Range("L140").Value = 88494.44500000001
$ws.Range("N140").Value = -98854.44500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1609.5
$ws.Range("I74").Value = 1146
$ws.Range("K74").Value = 1146
$ws.Range("M74").Value = -272
$ws.Range("H77").Value = 1609.5
$ws.Range("I77").Value = 1146
$ws.Range("K77").Value = 5730
$ws.Range("M77").Value = -1362
$ws.Range("H88").Value = 4981.95
$ws.Range("I88").Value = 1999.5
$ws.Range("J88").Value = 5313.3335
$ws.Range("K88").Value = 1999.5
$ws.Range("L88").Value = 5313.3335
$ws.Range("M88").Value = -1593.5
$ws.Range("N88").Value = -6125.3335
$ws.Range("H91").Value = 4981.95
$ws.Range("I91").Value = 1999.5
$ws.Range("J91").Value = 5313.3335
$ws.Range("K91").Value = 1999.5
$ws.Range("L91").Value = 5313.3335
$ws.Range("M91").Value = -595.5
$ws.Range("N91").Value = -8121.3335
$ws.Range("H122").Value = 5401.7144
$ws.Range("I122").Value = 5866.706
$ws.Range("K122").Value = 17600.118
$ws.Range("M122").Value = -15150.118
$ws.Range("H138").Value = 62660
$ws.Range("J138").Value = 62660
$ws.Range("L138").Value = 62660
$ws.Range("N138").Value = -72940
$ws.Range("H139").Value = 53238.332
$ws.Range("J139").Value = 53238.332
$ws.Range("L139").Value = 53238.332
$ws.Range("N139").Value = -63518.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2106.2856
$ws.Range("I94").Value = 1889.9
$ws.Range("J94").Value = 2647.25
$ws.Range("K94").Value = 1889.9
$ws.Range("L94").Value = 2647.25
$ws.Range("M94").Value = -1438.9
$ws.Range("N94").Value = -3549.25
$ws.Range("H107").Value = 3704839
$ws.Range("I107").Value = 4386936.5
$ws.Range("J107").Value = 2024
$ws.Range("K107").Value = 4386936.5
$ws.Range("L107").Value = 2024
$ws.Range("M107").Value = -4385016.5
$ws.Range("N107").Value = -5864
$ws.Range("H138").Value = 63226.668
$ws.Range("J138").Value = 63226.668
$ws.Range("L138").Value = 63226.668
$ws.Range("N138").Value = -73506.66800000001
$ws.Range("H140").Value = 89740
$ws.Range("J140").Value = 89740
$ws.Range("L140").Value = 89740
$ws.Range("N140").Value = -100100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4909.132
$ws.Range("I31").Value = 5530.8887
$ws.Range("K31").Value = 5530.8887
$ws.Range("M31").Value = -5235.8887
$ws.Range("H34").Value = 4909.132
$ws.Range("I34").Value = 5530.8887
$ws.Range("K34").Value = 5530.8887
$ws.Range("M34").Value = -5328.8887
$ws.Range("H68").Value = 17100.5
$ws.Range("J68").Value = 17100.5
$ws.Range("L68").Value = 17100.5
$ws.Range("N68").Value = -18598.5
$ws.Range("H71").Value = 17100.5
$ws.Range("J71").Value = 17100.5
$ws.Range("L71").Value = 51301.5
$ws.Range("N71").Value = -58789.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 1960
$ws.Range("J20").Value = 1960
$ws.Range("L20").Value = 5880
$ws.Range("N20").Value = -6334
$ws.Range("H56").Value = 2806316.5
$ws.Range("I56").Value = 2806316.5
$ws.Range("K56").Value = 2806316.5
$ws.Range("M56").Value = -2805786.5
$ws.Range("H70").Value = 1631.3636
$ws.Range("I70").Value = 870
$ws.Range("J70").Value = 1986.6666
$ws.Range("K70").Value = 2610
$ws.Range("L70").Value = 5959.9998
$ws.Range("M70").Value = -2295
$ws.Range("N70").Value = -6589.9998
$ws.Range("H73").Value = 1631.3636
$ws.Range("I73").Value = 870
$ws.Range("J73").Value = 1986.6666
$ws.Range("K73").Value = 2610
$ws.Range("L73").Value = 5959.9998
$ws.Range("M73").Value = -1518
$ws.Range("N73").Value = -8143.9998
$ws.Range("H75").Value = 626.6667
$ws.Range("I75").Value = 440
$ws.Range("J75").Value = 1000
$ws.Range("K75").Value = 1320
$ws.Range("L75").Value = 3000
$ws.Range("M75").Value = -322
$ws.Range("N75").Value = -4996
$ws.Range("H78").Value = 626.6667
$ws.Range("I78").Value = 440
$ws.Range("J78").Value = 1000
$ws.Range("K78").Value = 3960
$ws.Range("L78").Value = 9000
$ws.Range("M78").Value = 1032
$ws.Range("N78").Value = -18984
$ws.Range("H107").Value = 421.6
$ws.Range("I107").Value = 271.42856
$ws.Range("J107").Value = 553
$ws.Range("K107").Value = 814.28568
$ws.Range("L107").Value = 1659
$ws.Range("M107").Value = 1105.71432
$ws.Range("N107").Value = -5499
$ws.Range("H113").Value = 968.3299
$ws.Range("I113").Value = 616.6667
$ws.Range("J113").Value = 1017.97644
$ws.Range("K113").Value = 1850.0001
$ws.Range("L113").Value = 3053.92932
$ws.Range("M113").Value = 319.9999
$ws.Range("N113").Value = -7393.92932

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 38674.145
$ws.Range("I43").Value = 1350
$ws.Range("J43").Value = 53603.8
$ws.Range("K43").Value = 1350
$ws.Range("L43").Value = 53603.8
$ws.Range("M43").Value = -1199
$ws.Range("N43").Value = -53905.8
$ws.Range("H46").Value = 33142.855
$ws.Range("I46").Value = 27200
$ws.Range("J46").Value = 48000
$ws.Range("K46").Value = 27200
$ws.Range("L46").Value = 48000
$ws.Range("M46").Value = -27044
$ws.Range("N46").Value = -48312
$ws.Range("H57").Value = 18922.334
$ws.Range("I57").Value = 10000
$ws.Range("J57").Value = 20037.625
$ws.Range("K57").Value = 10000
$ws.Range("L57").Value = 20037.625
$ws.Range("M57").Value = -9180
$ws.Range("N57").Value = -21677.625
$ws.Range("H80").Value = 2268.3333
$ws.Range("I80").Value = 2268.3333
$ws.Range("K80").Value = 2268.3333
$ws.Range("M80").Value = -1270.3333
$ws.Range("H82").Value = 29975
$ws.Range("I82").Value = 10000
$ws.Range("J82").Value = 32828.57
$ws.Range("K82").Value = 10000
$ws.Range("L82").Value = 32828.57
$ws.Range("M82").Value = -9617
$ws.Range("N82").Value = -33594.57
$ws.Range("H83").Value = 2268.3333
$ws.Range("I83").Value = 2268.3333
$ws.Range("K83").Value = 11341.6665
$ws.Range("M83").Value = -6349.666499999999
$ws.Range("H85").Value = 29975
$ws.Range("I85").Value = 10000
$ws.Range("J85").Value = 32828.57
$ws.Range("K85").Value = 10000
$ws.Range("L85").Value = 32828.57
$ws.Range("M85").Value = -8674
$ws.Range("N85").Value = -35480.57
$ws.Range("H102").Value = 1904.591
$ws.Range("I102").Value = 1706
$ws.Range("J102").Value = 3162.3333
$ws.Range("K102").Value = 1706
$ws.Range("L102").Value = 3162.3333
$ws.Range("M102").Value = -84
$ws.Range("N102").Value = -6406.3333
$ws.Range("H140").Value = 88996.336
$ws.Range("J140").Value = 88996.336
$ws.Range("L140").Value = 88996.336
$ws.Range("N140").Value = -99356.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1619.4
$ws.Range("I55").Value = 2320
$ws.Range("K55").Value = 2320
$ws.Range("M55").Value = -2147
$ws.Range("H68").Value = 2116.0967
$ws.Range("I68").Value = 2056.4348
$ws.Range("J68").Value = 2287.625
$ws.Range("K68").Value = 2056.4348
$ws.Range("L68").Value = 2287.625
$ws.Range("M68").Value = -1307.4348
$ws.Range("N68").Value = -3785.625
$ws.Range("H71").Value = 2116.0967
$ws.Range("I71").Value = 2056.4348
$ws.Range("J71").Value = 2287.625
$ws.Range("K71").Value = 10282.174
$ws.Range("L71").Value = 11438.125
$ws.Range("M71").Value = -6538.173999999999
$ws.Range("N71").Value = -18926.125
$ws.Range("H136").Value = 7967.1816
$ws.Range("I136").Value = 6626
$ws.Range("K136").Value = 19878
$ws.Range("M136").Value = -17328
$ws.Range("H138").Value = 44021.418
$ws.Range("J138").Value = 44021.418
$ws.Range("L138").Value = 44021.418
$ws.Range("N138").Value = -54301.418
$ws.Range("H139").Value = 67700
$ws.Range("J139").Value = 67700
$ws.Range("L139").Value = 67700
$ws.Range("N139").Value = -77980

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 680.4167
$ws.Range("I107").Value = 653.375
$ws.Range("J107").Value = 734.5
$ws.Range("K107").Value = 1960.125
$ws.Range("L107").Value = 2203.5
$ws.Range("M107").Value = -40.125
$ws.Range("N107").Value = -6043.5
$ws.Range("H136").Value = 1007.2258
$ws.Range("I136").Value = 925.5
$ws.Range("K136").Value = 2776.5
$ws.Range("M136").Value = -226.5
$ws.Range("H138").Value = 62050
$ws.Range("J138").Value = 62050
$ws.Range("L138").Value = 62050
$ws.Range("N138").Value = -72330
$ws.Range("H141").Value = 69016.42999999999
$ws.Range("J141").Value = 69016.42999999999
$ws.Range("L141").Value = 69016.42999999999
$ws.Range("N141").Value = -79376.42999999999
